# Updated cryptos list on Tue May 23 19:36:31 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '27.217.14'
$ws.Range("E2").Value = '  +0.86%  '
$ws.Range("D3").Value = '1.852.21'
$ws.Range("E3").Value = '  +1.41%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.22'
$ws.Range("E5").Value = '  +0.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.51%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4622'
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3716'
$ws.Range("E8").Value = '  +0.30%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07275'
$ws.Range("E9").Value = '  -0.79%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8863'
$ws.Range("E10").Value = '  +1.35%  '
$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.00'
$ws.Range("E11").Value = '  +1.08%  '
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07796'
$ws.Range("E12").Value = '  -1.89%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.368'
$ws.Range("E13").Value = '  +0.51%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.767.82'
$ws.Range("E14").Value = '  -4.92%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.505'
$ws.Range("E15").Value = '  -0.85%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.13'
$ws.Range("E16").Value = '  -0.08%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.002'
$ws.Range("E17").Value = '  -0.51%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008910'
$ws.Range("E18").Value = '  +0.57%  '
$ws.Range("E19").Value = '  -0.56%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.69'
$ws.Range("E20").Value = '  -0.71%  '
$ws.Range("D21").Value = '27.239.06'
$ws.Range("E21").Value = '  -0.45%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.055'
$ws.Range("E22").Value = '  -0.92%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.49'
$ws.Range("E23").Value = '  -0.43%  '
$ws.Range("D24").Value = '2.081.51'
$ws.Range("E24").Value = '  -2.70%  '
$ws.Range("E25").Value = '  +5.42%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '152.06'
$ws.Range("E26").Value = '  -0.73%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.40'
$ws.Range("E27").Value = '  -0.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.044'
$ws.Range("E28").Value = '  +0.28%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '115.60'
$ws.Range("E29").Value = '  +0.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.055'
$ws.Range("E30").Value = '  -1.55%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08821'
$ws.Range("E31").Value = '  -0.89%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.124'
$ws.Range("E32").Value = '  +5.33%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7666'
$ws.Range("E33").Value = '  +5.19%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.172'
$ws.Range("E34").Value = '  +3.74%  '
$ws.Range("E35").Value = '  +1.33%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.757'
$ws.Range("E36").Value = '  +11.02%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.081'
$ws.Range("E37").Value = '  +1.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05228'
$ws.Range("E38").Value = '  +0.08%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01936'
$ws.Range("E39").Value = '  -0.94%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.932'
$ws.Range("E40").Value = '  -0.26%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.064'
$ws.Range("E41").Value = '  -0.80%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5094'
$ws.Range("E42").Value = '  -1.22%  '
$ws.Range("E43").Value = '  +0.10%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.388'
$ws.Range("E44").Value = '  +2.56%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4785'
$ws.Range("E45").Value = '  -1.30%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.31'
$ws.Range("E46").Value = '  +0.88%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.000'
$ws.Range("E47").Value = '  -0.55%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '102.92'
$ws.Range("E48").Value = '  +0.17%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.634'
$ws.Range("E49").Value = '  +0.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06206'
$ws.Range("E50").Value = '  +0.22%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '65.45'
$ws.Range("E51").Value = '  +0.63%  '
